$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> index map used below: D=4, L=12, M=13, N=14, O=15, P=16, Q=17, R=18, S=19, T=20

# Row 2
$ws.Cells.Item(2, 4).Value = 44776
$ws.Cells.Item(2, 13).Value = 50
$ws.Cells.Item(2, 14).Value = 10000
$ws.Cells.Item(2, 16).Value = 10000
$ws.Cells.Item(2, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(2, 19).Value = 556
$ws.Cells.Item(2, 20).Value = 18

# Row 3
$ws.Cells.Item(3, 4).Value = 44776
$ws.Cells.Item(3, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(3, 19).Value = 444
$ws.Cells.Item(3, 20).Value = 18

# Row 4
$ws.Cells.Item(4, 4).Value = 44272
$ws.Cells.Item(4, 14).Value = 9000
$ws.Cells.Item(4, 15).Value = 10000
$ws.Cells.Item(4, 16).Value = 9500
$ws.Cells.Item(4, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(4, 19).Value = 633
$ws.Cells.Item(4, 20).Value = 15

# Row 5
$ws.Cells.Item(5, 4).Value = 44272
$ws.Cells.Item(5, 12).Value = "Segunda"
$ws.Cells.Item(5, 14).Value = 8000
$ws.Cells.Item(5, 15).Value = 8000
$ws.Cells.Item(5, 16).Value = 8000
$ws.Cells.Item(5, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(5, 19).Value = 533
$ws.Cells.Item(5, 20).Value = 15

# Row 6
$ws.Cells.Item(6, 4).Value = 44363
$ws.Cells.Item(6, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(6, 19).Value = 633
$ws.Cells.Item(6, 20).Value = 15

# Row 7
$ws.Cells.Item(7, 4).Value = 44307
$ws.Cells.Item(7, 13).Value = 50
$ws.Cells.Item(7, 14).Value = 10000
$ws.Cells.Item(7, 16).Value = 10000
$ws.Cells.Item(7, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(7, 19).Value = 556
$ws.Cells.Item(7, 20).Value = 18

# Row 8
$ws.Cells.Item(8, 4).Value = 44307
$ws.Cells.Item(8, 12).Value = "Segunda"
$ws.Cells.Item(8, 13).Value = 50
$ws.Cells.Item(8, 14).Value = 8000
$ws.Cells.Item(8, 15).Value = 8000
$ws.Cells.Item(8, 16).Value = 8000
$ws.Cells.Item(8, 19).Value = 444

# Row 9
$ws.Cells.Item(9, 4).Value = 44299
$ws.Cells.Item(9, 13).Value = 100
$ws.Cells.Item(9, 15).Value = 11000
$ws.Cells.Item(9, 16).Value = 10500
$ws.Cells.Item(9, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(9, 18).Value = "Región del Maule"
$ws.Cells.Item(9, 19).Value = 583

# Row 10
$ws.Cells.Item(10, 4).Value = 44299
$ws.Cells.Item(10, 14).Value = 9000
$ws.Cells.Item(10, 15).Value = 9000
$ws.Cells.Item(10, 16).Value = 9000
$ws.Cells.Item(10, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(10, 18).Value = "Región del Maule"
$ws.Cells.Item(10, 19).Value = 500

# Row 11
$ws.Cells.Item(11, 4).Value = 44698
$ws.Cells.Item(11, 13).Value = 50
$ws.Cells.Item(11, 15).Value = 10000
$ws.Cells.Item(11, 16).Value = 10000
$ws.Cells.Item(11, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(11, 19).Value = 556

# Row 12
$ws.Cells.Item(12, 4).Value = 44358
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 100
$ws.Cells.Item(12, 14).Value = 11000
$ws.Cells.Item(12, 15).Value = 12000
$ws.Cells.Item(12, 16).Value = 11500
$ws.Cells.Item(12, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(12, 19).Value = 639

# Row 13
$ws.Cells.Item(13, 4).Value = 44316
$ws.Cells.Item(13, 13).Value = 100
$ws.Cells.Item(13, 14).Value = 9000
$ws.Cells.Item(13, 16).Value = 9500
$ws.Cells.Item(13, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(13, 19).Value = 528

# Row 14
$ws.Cells.Item(14, 4).Value = 44425
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 100
$ws.Cells.Item(14, 14).Value = 12000
$ws.Cells.Item(14, 15).Value = 13000
$ws.Cells.Item(14, 16).Value = 12500
$ws.Cells.Item(14, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(14, 19).Value = 694
